$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "OKTP"
$ws.Range("D2").Value = "Ok to Plate"
$ws.Range("O2").Value = "OKTP"
$ws.Range("P2").Value = "Ok to Plate`nIntegration Validation"
$ws.Rows(2).AutoFit()

$ws.Range("B3").Value = "Plate burn"
$ws.Range("D3").Value = "Plate Making"
$ws.Range("O3").Value = "Plate Making"
$ws.Range("P3").Value = "Plate Making"
